$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 16500
$ws.Range("J51").Value = 16500
$ws.Range("L51").Value = 16500
$ws.Range("N51").Value = -17468
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -11126
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -31632
$ws.Range("N72").Value = -62736
$ws.Range("H80").Value = 3075
$ws.Range("J80").Value = 3360
$ws.Range("L80").Value = 10080
$ws.Range("N80").Value = -12076
$ws.Range("H83").Value = 3075
$ws.Range("J83").Value = 3360
$ws.Range("L83").Value = 30240
$ws.Range("N83").Value = -40224
$ws.Range("H111").Value = 2310.9
$ws.Range("I111").Value = 2368.3333
$ws.Range("K111").Value = 7104.999899999999
$ws.Range("M111").Value = -4037.999899999999
$ws.Range("H135").Value = 4586.3335
$ws.Range("I135").Value = 4586.3335
$ws.Range("K135").Value = 41277.0015
$ws.Range("M135").Value = -38742.0015
$ws.Range("H137").Value = 3079.9092
$ws.Range("I137").Value = 1975.8
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 5927.4
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -3377.4
$ws.Range("N137").Value = -17100
$ws.Range("H138").Value = 4175.9
$ws.Range("J138").Value = 5049.2
$ws.Range("L138").Value = 15147.6
$ws.Range("N138").Value = -25427.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3093.1428
$ws.Range("I2").Value = 153.5
$ws.Range("K2").Value = 153.5
$ws.Range("M2").Value = -40.5
$ws.Range("H32").Value = 3671.9412
$ws.Range("J32").Value = 22222
$ws.Range("L32").Value = 22222
$ws.Range("N32").Value = -22796
$ws.Range("H97").Value = 224.46153
$ws.Range("I97").Value = 210
$ws.Range("K97").Value = 210
$ws.Range("M97").Value = 286
$ws.Range("H116").Value = 3093.1428
$ws.Range("I116").Value = 153.5
$ws.Range("K116").Value = 153.5
$ws.Range("M116").Value = 2140.5
$ws.Range("H122").Value = 2922.8333
$ws.Range("J122").Value = 2875
$ws.Range("L122").Value = 8625
$ws.Range("N122").Value = -13525
$ws.Range("H132").Value = 853
$ws.Range("I132").Value = 853
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2559
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -29
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3093.1428
$ws.Range("I3").Value = 153.5
$ws.Range("K3").Value = 153.5
$ws.Range("M3").Value = -39.5
$ws.Range("H20").Value = 4638.4614
$ws.Range("I20").Value = 1530.3
$ws.Range("K20").Value = 1530.3
$ws.Range("M20").Value = -1283.3
$ws.Range("H86").Value = 1682.7778
$ws.Range("I86").Value = 1476
$ws.Range("K86").Value = 1476
$ws.Range("M86").Value = -353
$ws.Range("H89").Value = 1682.7778
$ws.Range("I89").Value = 1476
$ws.Range("K89").Value = 7380
$ws.Range("M89").Value = -1764

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1329.3334
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1329.3334
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1329.3334
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1903.3334
$ws.Range("H22").Value = 1233.1666
$ws.Range("I22").Value = 750
$ws.Range("K22").Value = 750
$ws.Range("M22").Value = -400
$ws.Range("H58").Value = 9727.700000000001
$ws.Range("I58").Value = 9609.714
$ws.Range("J58").Value = 10003
$ws.Range("K58").Value = 9609.714
$ws.Range("L58").Value = 10003
$ws.Range("M58").Value = -9406.714
$ws.Range("N58").Value = -10409
$ws.Range("H99").Value = 5763.8667
$ws.Range("I99").Value = 5173.5
$ws.Range("J99").Value = 6944.6
$ws.Range("K99").Value = 5173.5
$ws.Range("L99").Value = 6944.6
$ws.Range("M99").Value = -3675.5
$ws.Range("N99").Value = -9940.6
$ws.Range("H113").Value = 1329.3334
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1329.3334
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1329.3334
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5669.3334
$ws.Range("H126").Value = 5763.8667
$ws.Range("I126").Value = 5173.5
$ws.Range("J126").Value = 6944.6
$ws.Range("K126").Value = 15520.5
$ws.Range("L126").Value = 20833.8
$ws.Range("M126").Value = -13050.5
$ws.Range("N126").Value = -25773.8
$ws.Range("H132").Value = 14842.714
$ws.Range("I132").Value = 13999
$ws.Range("K132").Value = 41997
$ws.Range("M132").Value = -39467
$ws.Range("H136").Value = 9727.700000000001
$ws.Range("I136").Value = 9609.714
$ws.Range("J136").Value = 10003
$ws.Range("K136").Value = 28829.142
$ws.Range("L136").Value = 30009
$ws.Range("M136").Value = -26279.142
$ws.Range("N136").Value = -35109

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1402
$ws.Range("J18").Value = 1825.8334
$ws.Range("L18").Value = 5477.5002
$ws.Range("N18").Value = -5815.5002
$ws.Range("H98").Value = 531.2308
$ws.Range("I98").Value = 535.5
$ws.Range("J98").Value = 529.3333
$ws.Range("K98").Value = 1606.5
$ws.Range("L98").Value = 1587.9999
$ws.Range("M98").Value = -108.5
$ws.Range("N98").Value = -4583.9999
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1299.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 1299.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 1299.5
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -4543.5
$ws.Range("H122").Value = 4571.2856
$ws.Range("I122").Value = 4571.2856
$ws.Range("K122").Value = 13713.8568
$ws.Range("M122").Value = -11263.8568
$ws.Range("H126").Value = 6249.4
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6249.4
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 18748.2
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -23688.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2958.5
$ws.Range("J46").Value = 2958.5
$ws.Range("L46").Value = 2958.5
$ws.Range("N46").Value = -3334.5
$ws.Range("H122").Value = 4999
$ws.Range("I122").Value = 4999
$ws.Range("K122").Value = 14997
$ws.Range("M122").Value = -12547
$ws.Range("H132").Value = 4284.923
$ws.Range("I132").Value = 3880.4
$ws.Range("J132").Value = 5633.3335
$ws.Range("K132").Value = 11641.2
$ws.Range("L132").Value = 16900.0005
$ws.Range("M132").Value = -9111.200000000001
$ws.Range("N132").Value = -21960.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1982.6666
$ws.Range("I96").Value = 1549.2858
$ws.Range("J96").Value = 3499.5
$ws.Range("K96").Value = 1549.2858
$ws.Range("L96").Value = 3499.5
$ws.Range("M96").Value = -176.2858000000001
$ws.Range("N96").Value = -6245.5
$ws.Range("H122").Value = 3999
$ws.Range("I122").Value = 3999
$ws.Range("K122").Value = 11997
$ws.Range("M122").Value = -9547
$ws.Range("H132").Value = 3474.9
$ws.Range("I132").Value = 3470.5715
$ws.Range("J132").Value = 3485
$ws.Range("K132").Value = 10411.7145
$ws.Range("L132").Value = 10455
$ws.Range("M132").Value = -7881.7145
$ws.Range("N132").Value = -15515
$ws.Range("H136").Value = 3942
$ws.Range("I136").Value = 3640
$ws.Range("K136").Value = 10920
$ws.Range("M136").Value = -8370
